# Apply edits to NEW_HAMPSHIRE_2018.xlsx as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Fix capitalization of "de"/"del"/"los" -> "De"/"Del"/"Los" in municipality/state names
$ws.Range("B2").Value = "Comitán De Domínguez"
$ws.Range("B12").Value = "Villa De Álvarez"
$ws.Range("A14").Value = "Ciudad De México"
$ws.Range("A24").Value = "Estado De México"
$ws.Range("B24").Value = "Ecatepec De Morelos"
$ws.Range("B28").Value = "Jaral Del Progreso"
$ws.Range("B30").Value = "Acapulco De Juárez"
$ws.Range("B31").Value = "Chilpancingo De Los Bravo"
$ws.Range("B32").Value = "Huitzuco De Los Figueroa"
$ws.Range("B34").Value = "Mártir De Cuilapan"
$ws.Range("B37").Value = "Taxco De Alarcón"
$ws.Range("B42").Value = "Autlán De Navarro"
$ws.Range("B50").Value = "Zapotitlán De Vadillo"
$ws.Range("B60").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B61").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B62").Value = "Putla Villa De Guerrero"
$ws.Range("B67").Value = "Palmar De Bravo"
$ws.Range("B70").Value = "Cadereyta De Montes"
$ws.Range("B71").Value = "Landa De Matamoros"
$ws.Range("B82").Value = "Juchique De Ferrer"

# 3. Remove trailing footer rows (90-94), which also removes the empty gap row 89
$ws.Range("A89:A94").EntireRow.Delete()
